$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = 4.9
$ws.Range("P3").Value = 1.72
$ws.Range("AI3").Value = 28
$ws.Range("AM3").Value = 55
$ws.Range("G6").Value = 2.3
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 2.82
$ws.Range("J6").Value = 2.87
$ws.Range("K6").Value = 2.12
$ws.Range("L6").Value = 3.35
$ws.Range("M6").Value = 1.33
$ws.Range("O6").Value = 1.98
$ws.Range("P6").Value = 1.65
$ws.Range("R6").Value = 1.25
$ws.Range("S6").Value = 1.38
$ws.Range("T6").Value = 2.57
$ws.Range("U6").Value = 1.78
$ws.Range("V6").Value = 1.82
$ws.Range("W6").Value = 7.4
$ws.Range("X6").Value = 10.75
$ws.Range("Y6").Value = 9.25
$ws.Range("Z6").Value = 22
$ws.Range("AA6").Value = 19.5
$ws.Range("AB6").Value = 32
$ws.Range("AD6").Value = 6.5
$ws.Range("AF6").Value = 75
$ws.Range("AH6").Value = 8.5
$ws.Range("AI6").Value = 13.5
$ws.Range("AJ6").Value = 10.5
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 25
$ws.Range("AM6").Value = 37
$ws.Range("G8").Value = 2.07
$ws.Range("H8").Value = 3.65
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 2.45
$ws.Range("K8").Value = 2.27
$ws.Range("M8").Value = 1.22
$ws.Range("N8").Value = 3.45
$ws.Range("O8").Value = 1.65
$ws.Range("P8").Value = 1.98
$ws.Range("Q8").Value = 2.57
$ws.Range("R8").Value = 1.39
$ws.Range("U8").Value = 1.6
$ws.Range("V8").Value = 2.07
$ws.Range("W8").Value = 9
$ws.Range("X8").Value = 10.75
$ws.Range("Y8").Value = 8.75
$ws.Range("Z8").Value = 19
$ws.Range("AA8").Value = 15.5
$ws.Range("AB8").Value = 24
$ws.Range("AC8").Value = 13
$ws.Range("AD8").Value = 7.2
$ws.Range("AE8").Value = 13.5
$ws.Range("AF8").Value = 50
$ws.Range("AG8").Value = 350
$ws.Range("AH8").Value = 11
$ws.Range("AI8").Value = 16.5
$ws.Range("AJ8").Value = 10.75
$ws.Range("AK8").Value = 37
$ws.Range("AL8").Value = 24
$ws.Range("AM8").Value = 29
$ws.Range("S10").Value = 1.41
$ws.Range("T10").Value = 2.62
$ws.Range("H11").Value = 3
$ws.Range("S11").Value = 1.53
$ws.Range("T11").Value = 2.38
$ws.Range("X11").Value = 10
$ws.Range("G12").Value = 3.2
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 2.38
$ws.Range("J12").Value = 4
$ws.Range("L12").Value = 3.25
$ws.Range("U12").Value = 2.1
$ws.Range("V12").Value = 1.67
$ws.Range("AC12").Value = 7
$ws.Range("AL12").Value = 23
$ws.Range("AN12").Value = 1.1
$ws.Range("AO12").Value = 7
$ws.Range("AP12").Value = 1.85
$ws.Range("AQ12").Value = 2
